$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.25
$ws.Range("K2").Value = 8.5
$ws.Range("N2").Value = 2.1
$ws.Range("O2").Value = 1.7
$ws.Range("Z2").Value = 8.5
$ws.Range("L3").Value = 1.25
$ws.Range("M3").Value = 3.75
$ws.Range("N3").Value = 1.8
$ws.Range("O3").Value = 2
$ws.Range("G4").Value = 2.25
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 1.11
$ws.Range("K4").Value = 6.5
$ws.Range("Z4").Value = 6.5
$ws.Range("AB4").Value = 17
$ws.Range("AH4").Value = 41
$ws.Range("G14").Value = 3.55
$ws.Range("I14").Value = 1.93
$ws.Range("Q14").Value = 2.35
$ws.Range("T14").Value = 7.4
$ws.Range("U14").Value = 14
$ws.Range("V14").Value = 10.75
$ws.Range("W14").Value = 37
$ws.Range("X14").Value = 29
$ws.Range("AE14").Value = 5.2
$ws.Range("AF14").Value = 7
$ws.Range("AG14").Value = 7.5
$ws.Range("AH14").Value = 13
$ws.Range("AI14").Value = 14
$ws.Range("G15").Value = 3.4
$ws.Range("I15").Value = 2.3
$ws.Range("P15").Value = 1.5
$ws.Range("Q15").Value = 2.5
$ws.Range("R15").Value = 2
$ws.Range("S15").Value = 1.75
$ws.Range("T15").Value = 8.5
$ws.Range("Z15").Value = 7
$ws.Range("AH15").Value = 21
$ws.Range("AI15").Value = 21
$ws.Range("AJ15").Value = 34
$ws.Range("N16").Value = 2.15
$ws.Range("O16").Value = 1.67
$ws.Range("I18").Value = 3.8
$ws.Range("K18").Value = 6.9
$ws.Range("Q18").Value = 2.57
$ws.Range("T18").Value = 6.7
$ws.Range("W18").Value = 15.5
$ws.Range("Z18").Value = 6.9
$ws.Range("AC18").Value = 80
$ws.Range("AE18").Value = 10.25
$ws.Range("AF18").Value = 20
$ws.Range("AG18").Value = 13.5
$ws.Range("AI18").Value = 40
$ws.Range("I22").Value = 1.85
$ws.Range("R22").Value = 1.75
$ws.Range("S22").Value = 2
$ws.Range("T22").Value = 12
$ws.Range("AD22").Value = 201
$ws.Range("H31").Value = 6
$ws.Range("I31").Value = 7
$ws.Range("R31").Value = 1.57
$ws.Range("S31").Value = 2.25
$ws.Range("U31").Value = 9.5
$ws.Range("V31").Value = 10
$ws.Range("W31").Value = 10
$ws.Range("AB31").Value = 19
$ws.Range("AD31").Value = 126
$ws.Range("AH31").Value = 81
$ws.Range("AJ31").Value = 41
$ws.Range("J32").Value = 1.05
$ws.Range("K32").Value = 8.5
$ws.Range("N32").Value = 1.93
$ws.Range("O32").Value = 1.88
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = 6.5
$ws.Range("I33").Value = 1.22
$ws.Range("K33").Value = 15
$ws.Range("R33").Value = 1.62
$ws.Range("S33").Value = 2.2
$ws.Range("AC33").Value = 41
$ws.Range("AD33").Value = 126
$ws.Range("AF33").Value = 9.5
$ws.Range("AH33").Value = 9.5
$ws.Range("N34").Value = 1.98
$ws.Range("O34").Value = 1.83
